$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: clear the "RL" marker in D3, replace with a single space
$ws.Range("D3").Value = " "

# New block describing "Download and Upload students in class (Class)"
# Row 62: main record row
$ws.Range("A62").Value = "Download and Upload students in class (Class)"
$ws.Range("B62").Value = 'name+"_StudData.txt"'
$ws.Range("C62").Value = "22CLC06_StudData.txt"
$ws.Range("D62").Value = "L1"
$ws.Range("E62").Value = "Number of students"
$ws.Range("F62").Value = "For loop"

# Row 63: continuation row
$ws.Range("D63").Value = "RL"
$ws.Range("E63").Value = "Student ID"
$ws.Range("F63").Value = "search"
$ws.Range("G63").Value = "Student"

# Update view: scroll position and selection to reflect the new content
$win = $excel.ActiveWindow
$win.ScrollRow = 52
$win.ScrollColumn = 1
$win.Zoom = 100
$ws.Range("A64").Select()
